$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker rows to append after the existing data (rows 454-458)
$newTickers = @("IMX-USD", "MNT-USD", "PEPE-USD", "GRT-USD", "TAO-USD")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newTickers[$i]
}
